# Fix "out of range" period rows: each class sheet (1반..4반) originally
# listed 8 periods (rows 2-9) with some invalid subjects (과학, 국어2, 수학2).
# Trim every sheet down to 5 valid periods (rows 2-6) and rewrite the
# Mon-Fri (B:F) subject grid so only the 5 valid subjects are used.

$wb = $excel.ActiveWorkbook

$data = @{
  "1반" = @{
    2 = @("체육","사회","국어1","수학1","사회")
    3 = @("영어","국어1","체육","영어","수학1")
    4 = @("사회","수학1","영어","체육","국어1")
    5 = @("수학1","영어","사회","국어1","체육")
    6 = @("국어1","체육","수학1","사회","영어")
  }
  "2반" = @{
    2 = @("영어","수학1","수학1","체육","체육")
    3 = @("수학1","영어","사회","사회","국어1")
    4 = @("체육","체육","체육","수학1","사회")
    5 = @("국어1","사회","영어","영어","영어")
    6 = @("사회","국어1","국어1","국어1","수학1")
  }
  "3반" = @{
    2 = @("국어1","국어1","사회","영어","영어")
    3 = @("사회","수학1","국어1","국어1","사회")
    4 = @("수학1","영어","수학1","사회","체육")
    5 = @("체육","체육","체육","체육","수학1")
    6 = @("영어","사회","영어","수학1","국어1")
  }
  "4반" = @{
    2 = @("수학1","영어","체육","국어1","수학1")
    3 = @("국어1","체육","영어","수학1","체육")
    4 = @("영어","사회","국어1","영어","영어")
    5 = @("사회","국어1","수학1","사회","국어1")
    6 = @("체육","수학1","사회","체육","사회")
  }
}

$cols = @("B","C","D","E","F")

foreach ($sheetName in $data.Keys) {
  $ws = $wb.Worksheets.Item($sheetName)

  # Remove the now out-of-range periods (old rows 7, 8, 9 = periods 6, 7, 8).
  $ws.Range("A7:F9").EntireRow.Delete()

  $rowsForSheet = $data[$sheetName]
  foreach ($r in $rowsForSheet.Keys) {
    $vals = $rowsForSheet[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
      $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
  }
}
